$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (Total) sheet.
#    We clone the "2021-Q4" sheet (it already carries the right column
#    layout / header style) and then overwrite its contents, which is
#    the most reliable way to reproduce the existing formatting.
# ------------------------------------------------------------------
$srcSheet   = $wb.Worksheets.Item(3)   # "2021-Q4"
$totalSheet = $wb.Worksheets.Item(4)   # "总计" (will end up after the new sheet)
$srcSheet.Copy($totalSheet, $null)

$q1 = $wb.Worksheets.Item(4)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Make sure the index column (A2:A8) keeps the same style as the header
# sheet we cloned from (the source sheet only had 6 data rows, this one
# needs 7, so extend the formatting down to row 8). Copy/paste ranges of
# matching size so the paste doesn't spill into extra rows.
$q1.Range("A2:A3").Copy() | Out-Null
$q1.Range("A7:A8").PasteSpecial(-4122) # xlPasteFormats

$q1Data = @(
    @("008866", "博时产业新趋势灵活配置混合A", "7.89", "86.47", "2.76", "0.2178", 10),
    @("501098", "建信科技创新 3 年封闭运作灵活配置混合", "3.61", "76.11", "2.38", "0.0859", 9),
    @("010665", "博时高端装备混合A", "0.77", "87.86", "4.56", "0.0351", 7),
    @("002595", "博时工业4.0主题股票", "0.61", "87.15", "4.27", "0.0260", 8),
    @("010666", "博时高端装备混合C", "0.41", "87.86", "4.56", "0.0187", 7),
    @("002567", "大成国家安全主题灵活配置混合", "0.34", "52.90", "3.63", "0.0123", 8),
    @("008867", "博时产业新趋势灵活配置混合C", "0.40", "86.47", "2.76", "0.0110", 10)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: add a new top row for "2022-Q1" with
#    (持有数量=7, 持有市值=0.41), pushing the previous rows down and
#    renumbering the index column.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item(5)   # "总计"
$total.Rows.Item(2).Insert()

# Re-apply the data-row formatting (the blank inserted row borrows the
# header row's bold style by default) from the row right below it.
$total.Range("A3:D3").Copy() | Out-Null
$total.Range("A2:D2").PasteSpecial(-4122) # xlPasteFormats

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.41

for ($r = 3; $r -le 5; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
